$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 7

# Numeric cells
$ws.Cells.Item($row, 1).Value = 131118183    # A7 Id
$ws.Cells.Item($row, 2).Value = 57881        # B7 Taxonsorteringsordning
$ws.Cells.Item($row, 5).Value = 100049       # E7 TaxonId
$ws.Cells.Item($row, 17).Value = 730393      # Q7 Ost
$ws.Cells.Item($row, 18).Value = 7123085     # R7 Nord
$ws.Cells.Item($row, 19).Value = 10          # S7 Noggrannhet

# Text (inline string) cells
$ws.Cells.Item($row, 4).Value = "NT"                     # D7 Rödlistade
$ws.Cells.Item($row, 6).Value = "Spillkråka"              # F7 Artnamn
$ws.Cells.Item($row, 7).Value = "Dryocopus martius"        # G7 Vetenskapligt namn
$ws.Cells.Item($row, 8).Value = "(Linnaeus, 1758)"         # H7 Auktor
$ws.Cells.Item($row, 16).Value = "Stor-Lidmyran, Vb"       # P7 Lokalnamn
$ws.Cells.Item($row, 20).Value = "Västerbotten"            # T7 Län
$ws.Cells.Item($row, 21).Value = "Vindeln"                 # U7 Kommun
$ws.Cells.Item($row, 22).Value = "Västerbotten"            # V7 Provins
$ws.Cells.Item($row, 23).Value = "Degerfors"                # W7 Socken
$ws.Cells.Item($row, 49).Value = "Elin Kannerby"           # AW7 Rapportör
$ws.Cells.Item($row, 50).Value = "Elin Kannerby"           # AX7 Observatörer

# Date-looking text cells -- force text format first so Excel does not
# auto-convert the string into a date serial number.
$ws.Cells.Item($row, 25).NumberFormat = "@"
$ws.Cells.Item($row, 25).Value = "2026-02-08"               # Y7 Startdatum
$ws.Cells.Item($row, 27).NumberFormat = "@"
$ws.Cells.Item($row, 27).Value = "2026-02-08"               # AA7 Slutdatum

# Boolean cells
$ws.Cells.Item($row, 30).Value = $false      # AD7 Ej återfunnen
$ws.Cells.Item($row, 31).Value = $false      # AE7 Osäker artbestämning
$ws.Cells.Item($row, 33).Value = $false      # AG7 Ospontan

# Empty (but present) cells - set a formatting property first so the
# cell is materialized even though its value is blank.
$ws.Cells.Item($row, 9).NumberFormat = "General"
$ws.Cells.Item($row, 9).Value = ""           # I7 Antal
$ws.Cells.Item($row, 46).NumberFormat = "General"
$ws.Cells.Item($row, 46).Value = ""          # AT7 Bestämningsår
$ws.Cells.Item($row, 51).NumberFormat = "General"
$ws.Cells.Item($row, 51).Value = ""          # AY7 Projektnamn
